$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3598.3333
$ws.Range("I33").Value = 89.14286
$ws.Range("K33").Value = 89.14286
$ws.Range("M33").Value = 139.85714
$ws.Range("H51").Value = 8823.632
$ws.Range("J51").Value = 9293.75
$ws.Range("L51").Value = 9293.75
$ws.Range("N51").Value = -10261.75
$ws.Range("H74").Value = 11732.889
$ws.Range("I74").Value = 10941.714
$ws.Range("J74").Value = 14502
$ws.Range("K74").Value = 10941.714
$ws.Range("L74").Value = 14502
$ws.Range("M74").Value = -10005.714
$ws.Range("N74").Value = -16374
$ws.Range("H77").Value = 11732.889
$ws.Range("I77").Value = 10941.714
$ws.Range("J77").Value = 14502
$ws.Range("K77").Value = 54708.57
$ws.Range("L77").Value = 72510
$ws.Range("M77").Value = -50028.57
$ws.Range("N77").Value = -81870
$ws.Range("H113").Value = 9370.909
$ws.Range("I113").Value = 5000
$ws.Range("J113").Value = 9808
$ws.Range("K113").Value = 5000
$ws.Range("L113").Value = 9808
$ws.Range("M113").Value = -1746
$ws.Range("N113").Value = -16316
$ws.Range("H132").Value = 1798
$ws.Range("I132").Value = 1798
$ws.Range("K132").Value = 5394
$ws.Range("M132").Value = -2864
$ws.Range("H137").Value = 4706.154
$ws.Range("I137").Value = 1020.25
$ws.Range("J137").Value = 6344.3335
$ws.Range("K137").Value = 3060.75
$ws.Range("L137").Value = 19033.0005
$ws.Range("M137").Value = -510.75
$ws.Range("N137").Value = -24133.0005
$ws.Range("H138").Value = 4928.5
$ws.Range("I138").Value = 3247.8
$ws.Range("J138").Value = 5978.9375
$ws.Range("K138").Value = 9743.400000000001
$ws.Range("L138").Value = 17936.8125
$ws.Range("M138").Value = -4603.400000000001
$ws.Range("N138").Value = -28216.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2515.389
$ws.Range("I32").Value = 2596
$ws.Range("K32").Value = 2596
$ws.Range("M32").Value = -2309
$ws.Range("H37").Value = 30000
$ws.Range("J37").Value = 30000
$ws.Range("L37").Value = 30000
$ws.Range("N37").Value = -30546
$ws.Range("H61").Value = 8743.799999999999
$ws.Range("I61").Value = 7542.636
$ws.Range("K61").Value = 7542.636
$ws.Range("M61").Value = -7330.636
$ws.Range("H63").Value = 5329.3335
$ws.Range("I63").Value = 1592.8
$ws.Range("K63").Value = 1592.8
$ws.Range("M63").Value = -906.8
$ws.Range("H66").Value = 5329.3335
$ws.Range("I66").Value = 1592.8
$ws.Range("K66").Value = 7964
$ws.Range("M66").Value = -4532
$ws.Range("H136").Value = 8743.799999999999
$ws.Range("I136").Value = 7542.636
$ws.Range("K136").Value = 22627.908
$ws.Range("M136").Value = -20077.908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2856.2354
$ws.Range("I20").Value = 1194.1111
$ws.Range("J20").Value = 4726.125
$ws.Range("K20").Value = 1194.1111
$ws.Range("L20").Value = 4726.125
$ws.Range("M20").Value = -947.1111000000001
$ws.Range("N20").Value = -5220.125
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H86").Value = 6291.3335
$ws.Range("I86").Value = 4629
$ws.Range("K86").Value = 4629
$ws.Range("M86").Value = -3506
$ws.Range("H89").Value = 6291.3335
$ws.Range("I89").Value = 4629
$ws.Range("K89").Value = 23145
$ws.Range("M89").Value = -17529
$ws.Range("H105").Value = 12515.759
$ws.Range("I105").Value = 11622.292
$ws.Range("J105").Value = 16804.4
$ws.Range("K105").Value = 11622.292
$ws.Range("L105").Value = 16804.4
$ws.Range("M105").Value = -9875.291999999999
$ws.Range("N105").Value = -20298.4
$ws.Range("H107").Value = 1441.2858
$ws.Range("I107").Value = 1280
$ws.Range("J107").Value = 1844.5
$ws.Range("K107").Value = 1280
$ws.Range("L107").Value = 1844.5
$ws.Range("M107").Value = 640
$ws.Range("N107").Value = -5684.5
$ws.Range("H134").Value = 3862.7646
$ws.Range("J134").Value = 8698.888999999999
$ws.Range("L134").Value = 26096.667
$ws.Range("N134").Value = -31166.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3264.2778
$ws.Range("I22").Value = 1562.25
$ws.Range("J22").Value = 4625.9
$ws.Range("K22").Value = 1562.25
$ws.Range("L22").Value = 4625.9
$ws.Range("M22").Value = -1212.25
$ws.Range("N22").Value = -5325.9
$ws.Range("H31").Value = 27327.775
$ws.Range("I31").Value = 3104.762
$ws.Range("J31").Value = 45495.035
$ws.Range("K31").Value = 3104.762
$ws.Range("L31").Value = 45495.035
$ws.Range("M31").Value = -2809.762
$ws.Range("N31").Value = -46085.035
$ws.Range("H34").Value = 27327.775
$ws.Range("I34").Value = 3104.762
$ws.Range("J34").Value = 45495.035
$ws.Range("K34").Value = 3104.762
$ws.Range("L34").Value = 45495.035
$ws.Range("M34").Value = -2902.762
$ws.Range("N34").Value = -45899.035
$ws.Range("H99").Value = 6488.875
$ws.Range("I99").Value = 6470.6665
$ws.Range("K99").Value = 6470.6665
$ws.Range("M99").Value = -4972.6665
$ws.Range("H105").Value = 4306.6943
$ws.Range("I105").Value = 3726.1
$ws.Range("K105").Value = 3726.1
$ws.Range("M105").Value = -1979.1
$ws.Range("H112").Value = 98750
$ws.Range("J112").Value = 98750
$ws.Range("L112").Value = 98750
$ws.Range("N112").Value = -101704
$ws.Range("H126").Value = 6488.875
$ws.Range("I126").Value = 6470.6665
$ws.Range("K126").Value = 19411.9995
$ws.Range("M126").Value = -16941.9995
$ws.Range("H132").Value = 3024.6875
$ws.Range("I132").Value = 2579.3044
$ws.Range("J132").Value = 4162.8887
$ws.Range("K132").Value = 7737.9132
$ws.Range("L132").Value = 12488.6661
$ws.Range("M132").Value = -5207.9132
$ws.Range("N132").Value = -17548.6661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5135.6895
$ws.Range("J68").Value = 5796.88
$ws.Range("L68").Value = 17390.64
$ws.Range("N68").Value = -19012.64
$ws.Range("H71").Value = 5135.6895
$ws.Range("J71").Value = 5796.88
$ws.Range("L71").Value = 52171.92
$ws.Range("N71").Value = -60283.92
$ws.Range("H107").Value = 1257.091
$ws.Range("I107").Value = 1294.4375
$ws.Range("J107").Value = 1157.5
$ws.Range("K107").Value = 3883.3125
$ws.Range("L107").Value = 3472.5
$ws.Range("M107").Value = -1963.3125
$ws.Range("N107").Value = -7312.5
$ws.Range("H108").Value = 8533.166999999999
$ws.Range("I108").Value = 4041.3333
$ws.Range("J108").Value = 13025
$ws.Range("K108").Value = 12123.9999
$ws.Range("L108").Value = 39075
$ws.Range("M108").Value = -9243.999899999999
$ws.Range("N108").Value = -44835

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 34613
$ws.Range("J93").Value = 35890
$ws.Range("L93").Value = 35890
$ws.Range("N93").Value = -39634
$ws.Range("H113").Value = 4311.24
$ws.Range("I113").Value = 1400.3846
$ws.Range("J113").Value = 7464.6665
$ws.Range("K113").Value = 1400.3846
$ws.Range("L113").Value = 7464.6665
$ws.Range("M113").Value = 769.6153999999999
$ws.Range("N113").Value = -11804.6665
$ws.Range("H128").Value = 70642.22
$ws.Range("J128").Value = 70642.22
$ws.Range("L128").Value = 70642.22
$ws.Range("N128").Value = -80602.22

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3083.2144
$ws.Range("I7").Value = 3090.3845
$ws.Range("K7").Value = 3090.3845
$ws.Range("M7").Value = -2978.3845
$ws.Range("H16").Value = 2671
$ws.Range("I16").Value = 2449.1904
$ws.Range("K16").Value = 2449.1904
$ws.Range("M16").Value = -2279.1904
$ws.Range("H22").Value = 4161.875
$ws.Range("J22").Value = 4949.1665
$ws.Range("L22").Value = 4949.1665
$ws.Range("N22").Value = -5539.1665
$ws.Range("H27").Value = 4161.875
$ws.Range("J27").Value = 4949.1665
$ws.Range("L27").Value = 4949.1665
$ws.Range("N27").Value = -5163.1665
$ws.Range("H46").Value = 4326.375
$ws.Range("I46").Value = 1219
$ws.Range("J46").Value = 4770.2856
$ws.Range("K46").Value = 1219
$ws.Range("L46").Value = 4770.2856
$ws.Range("M46").Value = -1031
$ws.Range("N46").Value = -5146.2856
$ws.Range("H55").Value = 16667017
$ws.Range("I55").Value = 50000000
$ws.Range("J55").Value = 525
$ws.Range("K55").Value = 50000000
$ws.Range("L55").Value = 525
$ws.Range("M55").Value = -49999827
$ws.Range("N55").Value = -871
$ws.Range("H93").Value = 9916.673000000001
$ws.Range("I93").Value = 6001.1943
$ws.Range("K93").Value = 6001.1943
$ws.Range("M93").Value = -4753.1943
$ws.Range("H119").Value = 68210.5
$ws.Range("J119").Value = 68210.5
$ws.Range("L119").Value = 68210.5
$ws.Range("N119").Value = -77886.5
$ws.Range("H126").Value = 3083.2144
$ws.Range("I126").Value = 3090.3845
$ws.Range("K126").Value = 9271.1535
$ws.Range("M126").Value = -6801.1535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1085.3846
$ws.Range("I107").Value = 1168
$ws.Range("K107").Value = 3504
$ws.Range("M107").Value = -1584
$ws.Range("H126").Value = 3697.353
$ws.Range("I126").Value = 2818.1765
$ws.Range("J126").Value = 4576.5293
$ws.Range("K126").Value = 8454.529500000001
$ws.Range("L126").Value = 13729.5879
$ws.Range("M126").Value = -5984.529500000001
$ws.Range("N126").Value = -18669.5879
$ws.Range("H129").Value = 54619
$ws.Range("J129").Value = 54619
$ws.Range("L129").Value = 54619
$ws.Range("N129").Value = -64619
